# Update workbook for "Add data for 2021-09-23":
#  - rename sheet / update "through" date from Sept 14 to Sept 15
#  - update the corresponding header label in the data grid
#  - bump a handful of existing monthly counts
#  - add a handful of brand-new (previously empty) monthly counts

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet name ---
$ws.Name = "Through 2021-09-15"

# --- Header label for the "through" column (B1 / shared string) ---
$ws.Range("B1").Value = "September 2021 (through September 15)"

# --- Updates to existing cell values ---
$ws.Range("K3").Value = 7    # North Lawndale,   September 2020: 6 -> 7
$ws.Range("B4").Value = 3    # Humboldt Park,     September 2021: 2 -> 3
$ws.Range("B5").Value = 6    # Austin,            September 2021: 5 -> 6
$ws.Range("K5").Value = 4    # Austin,            September 2020: 3 -> 4
$ws.Range("AC5").Value = 5   # Austin,            September 2018: 4 -> 5
$ws.Range("AL5").Value = 3   # Austin,            September 2017: 2 -> 3
$ws.Range("AC6").Value = 2   # Roseland,          September 2018: 1 -> 2
$ws.Range("B13").Value = 5   # Chatham,           September 2021: 4 -> 5
$ws.Range("T36").Value = 2   # West Elsdon,       September 2019: 1 -> 2
$ws.Range("K55").Value = 4   # Grand Crossing,    September 2020: 3 -> 4

# --- New (previously empty) cell values ---
$ws.Range("AU11").Value = 1  # Little Italy, UIC, September 2016
$ws.Range("B12").Value = 1   # Avondale,          September 2021
$ws.Range("AU12").Value = 1  # Avondale,          September 2016
$ws.Range("T17").Value = 1   # Grand Boulevard,   September 2019
$ws.Range("AL46").Value = 1  # South Deering,     September 2017
$ws.Range("B48").Value = 1   # Old Town,          September 2021
$ws.Range("B54").Value = 1   # Hermosa,           September 2021
$ws.Range("K60").Value = 1   # Armour Square,     September 2020
